$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45179
# (2023-09-10) to 45180 (2023-09-11) for every data row (rows 2-308).
$ws.Range("C2:C308").Value = 45180
